$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Class Statistics (K/L columns near top of sheet) ---
# Missing Sessions: 3 -> 9
$ws.Range("L7").Value = 9
# Pending Sessions: 114 -> 108
$ws.Range("L8").Value = 108

# --- Swap "Recorded By" name order for many session rows:
#     "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System" ---
$gRows = @(8, 9, 10, 12, 14, 15, 17, 18, 34, 35, 36, 38, 40, 41, 43, 44, 60, 61, 62, 64, 66, 67, 69, 70, 86, 87, 88, 90, 92, 93, 95, 96, 112, 113, 114, 116, 118, 119, 121, 122, 138, 139, 140, 142, 144, 145, 147, 148, 164, 167, 170, 174, 191, 194, 197, 201, 218, 221, 224, 228, 245, 248, 251, 255, 272, 275, 278, 282, 299, 302, 305, 309)
foreach ($r in $gRows) {
    $ws.Range("G$r").Value = "dnasr281@gmail.com, System"
}

# --- Group summary table (K:S columns, rows 21-26): update Not Recorded (P) and
#     Pending (Q) session counts for the B1D1/B1D2/B1E1/B1E2/B1F1/B1F2 groups ---
$ws.Range("P21").Value = 1
$ws.Range("Q21").Value = 9

$ws.Range("P22").Value = 1
$ws.Range("Q22").Value = 9

$ws.Range("P23").Value = 1
$ws.Range("Q23").Value = 9

$ws.Range("P24").Value = 2
$ws.Range("Q24").Value = 9

$ws.Range("P25").Value = 1
$ws.Range("Q25").Value = 9

$ws.Range("P26").Value = 1
$ws.Range("Q26").Value = 9

# --- Rows whose final (most-recent) session flips from "Pending" to "Not Recorded":
#     restyle A:I to the "Not Recorded" look (style of row 132, a known
#     Not-Recorded template row) and update the Status text in column I ---
$templateRange = $ws.Range("A132:I132")
$notRecordedRows = @(175, 202, 229, 256, 283, 310)
foreach ($r in $notRecordedRows) {
    $templateRange.Copy()
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122)
    $ws.Range("I$r").Value = "Not Recorded"
}
